# Replace old annotation with annotation of improved bins
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: rename "new_AP_ID" (H18) to "new_GOLD_ID"
$ws.Cells.Item(18, 8).Value = "new_GOLD_ID"

# Add new GOLD ID numeric values in column H for rows 19-29
$ws.Cells.Item(19, 8).Value = 2757320395
$ws.Cells.Item(20, 8).Value = 2757320396
$ws.Cells.Item(21, 8).Value = 2757320397
$ws.Cells.Item(22, 8).Value = 2757320398
$ws.Cells.Item(23, 8).Value = 2757320399
$ws.Cells.Item(24, 8).Value = 2757320400
$ws.Cells.Item(25, 8).Value = 2757320401
$ws.Cells.Item(26, 8).Value = 2757320402
$ws.Cells.Item(27, 8).Value = 2757320405
$ws.Cells.Item(28, 8).Value = 2757320403
$ws.Cells.Item(29, 8).Value = 2757320404

# Update column G (new_bin) for rows 23 and 24 with the refined bin names
$ws.Cells.Item(23, 7).Value = "B72-73_Su13.BD.MM15.SN.C_rebin5-6_refined1"
$ws.Cells.Item(24, 7).Value = "B72-73_Su13.BD.MM15.SN.C_rebin5-6_refined2"

# Update column I (new_IMG_ID) values with the new annotation for improved bins
$ws.Cells.Item(19, 9).Value = "Limnohabitans sp. bin L8r"
$ws.Cells.Item(20, 9).Value = "Limnohabitans sp. bin L5r"
$ws.Cells.Item(21, 9).Value = "Limnohabitans sp. bin L6r"
$ws.Cells.Item(22, 9).Value = "Limnohabitans sp. bin L7r"
$ws.Cells.Item(23, 9).Value = "Limnohabitans sp. bin L3r"
$ws.Cells.Item(24, 9).Value = "Limnohabitans sp. bin L1r"
$ws.Cells.Item(25, 9).Value = "Limnohabitans sp. bin L02"
$ws.Cells.Item(26, 9).Value = "Limnohabitans sp. bin L01"
$ws.Cells.Item(27, 9).Value = "Unclassified Betaproteobacteria bin B4r"
$ws.Cells.Item(28, 9).Value = "Limnohabitans sp. bin L03"
$ws.Cells.Item(29, 9).Value = "Limnohabitans sp. bin L2r"

# Row 30 no longer has new_IMG_ID / NCBI ID values - clear columns I and J
$ws.Cells.Item(30, 9).ClearContents()
$ws.Cells.Item(30, 10).ClearContents()

# Update the saved view state (scroll position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I26").Select()
